$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new titration record for 11/8/2019 (new CRM, fresh pH buffers, new acid dosing hose)
# Copy the date-format style from the prior row so A62 matches the existing date column formatting
$ws.Range("A61").Copy()
$ws.Range("A62").PasteSpecial(-4122)
$ws.Range("A62").Value = 43777

$ws.Range("B62").Value = 2142.32205763918
$ws.Range("C62").Value = 2207.0300000000002
$ws.Range("D62").Formula = "=100*(B62-C62)/C62"
$ws.Range("E62").Value = 169
$ws.Range("F62").Value = "new crm, fresh pH cal solutions, new hose for acid dosing"

# Move/extend the visible selection the way the author left it
$ws.Range("F63").Select()
